{"js": "// Apply the text replacements described by the diff.\nconst replacements = [\n  [\n    \"Play Going Underground Free: Exciting Slot Game Review\",\n    \"Play Going Underground Slot for Free\",\n  ],\n  [\n    \"Dynamic and exciting gameplay with high chances of winning payouts\",\n    \"Dynamic and exciting gameplay\",\n  ],\n  [\n    \"Impeccable graphics with interactive quality that enhances the gameplay experience\",\n    \"Impeccably designed graphics\",\n  ],\n  [\n    \"Multiple random modifiers and bonus rounds that keep the game continuously entertaining\",\n    \"Numerous bonus features and rounds\",\n  ],\n  [\n    \"Distinct and immersive theme that transports players back in time to explore the wonders of underground London\",\n    \"Distinct and immersive theme\",\n  ],\n  [\n    \"Limited number of paylines may not appeal to some players\",\n    \"Limited number of paylines\",\n  ],\n  [\n    \"No progressive jackpot feature\",\n    \"High volatility\",\n  ],\n  [\n    \"Experience the wonders of London's underground with Going Underground free slot game review. Discover the distinct and immersive theme and bonus features today!\",\n    \"Read our review of Going Underground, an exciting online slot game. Play for free and experience dynamic gameplay and numerous bonus features.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"Play Going Underground Free: Exciting Slot Game Review\"; New = \"Play Going Underground Slot for Free\" },\n    @{ Old = \"Dynamic and exciting gameplay with high chances of winning payouts\"; New = \"Dynamic and exciting gameplay\" },\n    @{ Old = \"Impeccable graphics with interactive quality that enhances the gameplay experience\"; New = \"Impeccably designed graphics\" },\n    @{ Old = \"Multiple random modifiers and bonus rounds that keep the game continuously entertaining\"; New = \"Numerous bonus features and rounds\" },\n    @{ Old = \"Distinct and immersive theme that transports players back in time to explore the wonders of underground London\"; New = \"Distinct and immersive theme\" },\n    @{ Old = \"Limited number of paylines may not appeal to some players\"; New = \"Limited number of paylines\" },\n    @{ Old = \"No progressive jackpot feature\"; New = \"High volatility\" },\n    @{ Old = \"Experience the wonders of London's underground with Going Underground free slot game review. Discover the distinct and immersive theme and bonus features today!\"; New = \"Read our review of Going Underground, an exciting online slot game. Play for free and experience dynamic gameplay and numerous bonus features.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
